$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 296.6
$ws.Range("I18").Value = 296.6
$ws.Range("K18").Value = 296.6
$ws.Range("M18").Value = -12.60000000000002
$ws.Range("H70").Value = 2986.6365
$ws.Range("J70").Value = 2255.0833
$ws.Range("L70").Value = 6765.249899999999
$ws.Range("N70").Value = -7305.249899999999
$ws.Range("H73").Value = 2986.6365
$ws.Range("J73").Value = 2255.0833
$ws.Range("L73").Value = 6765.249899999999
$ws.Range("N73").Value = -8637.249899999999
$ws.Range("H95").Value = 38491.5
$ws.Range("J95").Value = 35199.8
$ws.Range("L95").Value = 35199.8
$ws.Range("N95").Value = -40691.8
$ws.Range("H107").Value = 1652.1111
$ws.Range("J107").Value = 1450
$ws.Range("L107").Value = 1450
$ws.Range("N107").Value = -5290
$ws.Range("H113").Value = 16807.428
$ws.Range("I113").Value = 29417.334
$ws.Range("J113").Value = 7350
$ws.Range("K113").Value = 29417.334
$ws.Range("L113").Value = 7350
$ws.Range("M113").Value = -26163.334
$ws.Range("N113").Value = -13858
$ws.Range("H138").Value = 1836.1718
$ws.Range("I138").Value = 1178.1
$ws.Range("J138").Value = 2122.2898
$ws.Range("K138").Value = 3534.3
$ws.Range("L138").Value = 6366.8694
$ws.Range("M138").Value = 1605.7
$ws.Range("N138").Value = -16646.8694

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 7189.3335
$ws.Range("I22").Value = 775
$ws.Range("K22").Value = 775
$ws.Range("M22").Value = -476
$ws.Range("H32").Value = 3291.0173
$ws.Range("I32").Value = 2533.3062
$ws.Range("J32").Value = 7416.3335
$ws.Range("K32").Value = 2533.3062
$ws.Range("L32").Value = 7416.3335
$ws.Range("M32").Value = -2246.3062
$ws.Range("N32").Value = -7990.3335
$ws.Range("H45").Value = 3822.1875
$ws.Range("I45").Value = 4220.7
$ws.Range("K45").Value = 4220.7
$ws.Range("M45").Value = -3843.7
$ws.Range("H74").Value = 12829.723
$ws.Range("I74").Value = 1718.2667
$ws.Range("J74").Value = 68387
$ws.Range("K74").Value = 1718.2667
$ws.Range("L74").Value = 68387
$ws.Range("M74").Value = -844.2666999999999
$ws.Range("N74").Value = -70135
$ws.Range("H77").Value = 12829.723
$ws.Range("I77").Value = 1718.2667
$ws.Range("J77").Value = 68387
$ws.Range("K77").Value = 8591.333499999999
$ws.Range("L77").Value = 341935
$ws.Range("M77").Value = -4223.333499999999
$ws.Range("N77").Value = -350671
$ws.Range("H95").Value = 9750
$ws.Range("J95").Value = 9750
$ws.Range("L95").Value = 9750
$ws.Range("N95").Value = -15242
$ws.Range("H102").Value = 3280.2
$ws.Range("I102").Value = 1978
$ws.Range("K102").Value = 1978
$ws.Range("M102").Value = -356
$ws.Range("H110").Value = 6949.2354
$ws.Range("I110").Value = 9957.182000000001
$ws.Range("J110").Value = 1434.6666
$ws.Range("K110").Value = 9957.182000000001
$ws.Range("L110").Value = 1434.6666
$ws.Range("M110").Value = -7912.182000000001
$ws.Range("N110").Value = -5524.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 49117.85
$ws.Range("I134").Value = 48667.61
$ws.Range("K134").Value = 146002.83
$ws.Range("M134").Value = -143467.83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 337010.44
$ws.Range("J4").Value = 375386.75
$ws.Range("L4").Value = 375386.75
$ws.Range("N4").Value = -375610.75
$ws.Range("H58").Value = 12934.609
$ws.Range("I58").Value = 4783.8438
$ws.Range("K58").Value = 4783.8438
$ws.Range("M58").Value = -4580.8438
$ws.Range("H63").Value = 15271
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 15271
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H134").Value = 31256318
$ws.Range("I134").Value = 2054.913
$ws.Range("J134").Value = 111128330
$ws.Range("K134").Value = 6164.739
$ws.Range("L134").Value = 333384990
$ws.Range("M134").Value = -3629.739
$ws.Range("N134").Value = -333390060
$ws.Range("H136").Value = 12934.609
$ws.Range("I136").Value = 4783.8438
$ws.Range("K136").Value = 14351.5314
$ws.Range("M136").Value = -11801.5314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 11.333333
$ws.Range("I16").Value = 9.5
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 28.5
$ws.Range("L16").Value = 45
$ws.Range("M16").Value = 144.5
$ws.Range("N16").Value = -391
$ws.Range("H20").Value = 399
$ws.Range("I20").Value = 399
$ws.Range("K20").Value = 1197
$ws.Range("M20").Value = -970
$ws.Range("H109").Value = 4763533
$ws.Range("I109").Value = 1899.6666
$ws.Range("K109").Value = 5698.9998
$ws.Range("M109").Value = -4658.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H102").Value = 3910135.2
$ws.Range("I102").Value = 6949885.5
$ws.Range("K102").Value = 6949885.5
$ws.Range("M102").Value = -6948263.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3058.7334
$ws.Range("I46").Value = 2360.75
$ws.Range("K46").Value = 2360.75
$ws.Range("M46").Value = -2172.75
$ws.Range("H132").Value = 3657000.8
$ws.Range("I132").Value = 1999.8
$ws.Range("J132").Value = 6702834.5
$ws.Range("K132").Value = 5999.4
$ws.Range("L132").Value = 20108503.5
$ws.Range("M132").Value = -3469.4
$ws.Range("N132").Value = -20113563.5
$ws.Range("H136").Value = 132397.23
$ws.Range("I136").Value = 669333
$ws.Range("K136").Value = 2007999
$ws.Range("M136").Value = -2005449

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1619.1818
$ws.Range("J4").Value = 1839
$ws.Range("L4").Value = 1839
$ws.Range("N4").Value = -2065
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H132").Value = 44112
$ws.Range("I132").Value = 16334
$ws.Range("J132").Value = 99668
$ws.Range("K132").Value = 49002
$ws.Range("L132").Value = 299004
$ws.Range("M132").Value = -46472
$ws.Range("N132").Value = -304064
